# update data collection model
#
# 1. On the "DataCollectionItem" sheet, remove the "sdtmAnnotation" column (Q1).
# 2. On the "SDTMTarget" sheet, insert a new first column "sdtmAnnotation",
#    shifting "sdtmVariable" and "sdtmTargetMapping" one column to the right.

$wb = $excel.ActiveWorkbook

$itemSheet = $wb.Worksheets.Item("DataCollectionItem")
$itemSheet.Range("Q1").ClearContents() | Out-Null

$targetSheet = $wb.Worksheets.Item("SDTMTarget")
$targetSheet.Range("A1").EntireColumn.Insert() | Out-Null
$targetSheet.Range("A1").Value = "sdtmAnnotation"
